# Revert "Merge branch 'pdf-cover-sheet' into 'master'" (MR !4389)
#
# Removes the PDF cover-sheet content that was previously inserted right
# after the title paragraph ({{ caseType }} / {{ caseId }}):
#   - the extra empty Title-styled paragraph
#   - the big 3-column/14-row header table (address/plots/applicant/tag/
#     municipality/authority/responsible/inputDate/description rows)
#   - the trailing paragraph holding the page-break run
# and the now-unused "TableHeading" paragraph style from styles.xml.

$d = $word.ActiveDocument

# --- locate the boundaries -------------------------------------------------

# End of the first paragraph ({{ caseType }}<br/>{{ caseId }}), i.e. right
# after its paragraph mark. Everything from here up to (but not including)
# the "{%p for section in sections %}" paragraph is the cover-sheet content
# that needs to go away.
$titleEnd = $d.Paragraphs.Item(1).Range.End

$afterRange = $d.Content
[void]$afterRange.Find.Execute("{%p for section in sections %}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$keepStart = $afterRange.Start

# --- remove the cover-sheet table ------------------------------------------
# (Ranges that span across a table boundary cannot be deleted in one shot,
# so the table itself is removed via Table.Delete(), and the plain
# paragraph ranges before/after it are deleted separately.)

$tbl = $d.Tables.Item(1)
$tblStart = $tbl.Range.Start
$tblEnd = $tbl.Range.End

# Delete the trailing paragraphs after the table (empty colored paragraph +
# the page-break paragraph) up to the start of the paragraph we keep.
if ($keepStart -gt $tblEnd) {
    $trailing = $d.Range($tblEnd, $keepStart)
    $trailing.Delete()
}

# Delete the table itself.
$tbl2 = $d.Tables.Item(1)
$tbl2.Delete()

# Delete the leading empty Title paragraph between the kept title and the
# (now deleted) table.
if ($tblStart -gt $titleEnd) {
    $leading = $d.Range($titleEnd, $tblStart)
    $leading.Delete()
}

# --- remove the now-unused "TableHeading" paragraph style ------------------

$style = $d.Styles("TableHeading")
$style.Delete()
